# Add the intro paragraphs and the "7 Screen Reader Compatibility" heading
# right after the "Write Up" title paragraph, before the existing blank
# paragraphs.

$d = $word.ActiveDocument

# Locate the title paragraph ("Write Up") and insert a new paragraph after it.
$titlePara = $d.Paragraphs.Item(1)
$rng = $titlePara.Range
$rng.Collapse(0)              # wdCollapseEnd
$rng.InsertParagraphAfter()

# Paragraph 2: short intro blurb.
$p2 = $d.Paragraphs.Item(2)
$p2.Style = "Normal"
$p2.Range.Text = "This week, will be a very short tutorials on screen readers with Python and Tkinter. Although screen readers are not something that tkinter would explicitly deal with, there are certainly a few things that can be done to enhance your program to work with screen readers more effectively, and we will be discussing those options in this tutorial. "

$rng2 = $p2.Range
$rng2.Collapse(0)
$rng2.InsertParagraphAfter()

# Paragraph 3: call-to-action sentence introducing the tutorial title.
$p3 = $d.Paragraphs.Item(3)
$p3.Style = "Normal"
$p3.Range.Text = "So, if you would like to learn a bit more about how you can aide screen readers in their endeavor to enable the handicapped to gain access to information too, then please join us for our new tutorial this week entitled:"

$rng3 = $p3.Range
$rng3.Collapse(0)
$rng3.InsertParagraphAfter()

# Paragraph 4: the new article heading.
$p4 = $d.Paragraphs.Item(4)
$p4.Style = "Heading1"
$p4.Range.Text = "7 Screen Reader Compatibility"
